# Update the 2025 (row 8) recurrence metrics with the latest figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 1017
$ws.Range("D8").Value = 166
$ws.Range("E8").Value = 851
$ws.Range("F8").Value = 6.808859721082855
$ws.Range("G8").Value = 83.67748279252703
$ws.Range("H8").Value = 16.32251720747296
